$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New layout: swap the "soft" (E) and "rigid" (F) columns - including the
# header labels in row 1 and all the data counts below them.
$lastRow = 21

for ($r = 1; $r -le $lastRow; $r++) {
    $eCell = $ws.Cells.Item($r, 5)   # column E
    $fCell = $ws.Cells.Item($r, 6)   # column F

    $eVal = $eCell.Value2
    $fVal = $fCell.Value2

    $eCell.Value2 = $fVal
    $fCell.Value2 = $eVal
}
